$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from H1 to the new header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New data values for columns I and J
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 5

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 4

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 4

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 3

$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 2
